$wb = $excel.ActiveWorkbook

# --- Add a new trailing worksheet ("Sheet1") after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add($null, $lastSheet)

# --- Populate the "Array Sets" sheet with the new row-numbering column (A) ---
$ws = $wb.Worksheets.Item("Array Sets")

$ws.Range("A4").Value  = 1
$ws.Range("A6").Value  = 2
$ws.Range("A8").Value  = 3
$ws.Range("A10").Value = 4
$ws.Range("A12").Value = 5
$ws.Range("A14").Value = 6
$ws.Range("A16").Value = 7
$ws.Range("A18").Value = 8
$ws.Range("A20").Value = 9
$ws.Range("A22").Value = 10

# New row 24 with an extra "dummy" entry
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "dummy"

# --- Make "Array Sets" the active sheet/tab, with the view scrolled & the
#     selection parked on the row just below the new data ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B25").Select() | Out-Null
